$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "69-20=49"
$t.Cell(1,2).Range.Text = "26+35=61"
$t.Cell(1,3).Range.Text = "9+30=39"
$t.Cell(1,4).Range.Text = "93-3=90"
$t.Cell(1,5).Range.Text = "54-38=16"
$t.Cell(2,1).Range.Text = "59+9=68"
$t.Cell(2,2).Range.Text = "82-52=30"
$t.Cell(2,3).Range.Text = "91-48=43"
$t.Cell(2,4).Range.Text = "26+14=40"
$t.Cell(2,5).Range.Text = "16+31=47"
$t.Cell(3,1).Range.Text = "19+9=28"
$t.Cell(3,2).Range.Text = "95-58=37"
$t.Cell(3,3).Range.Text = "31-14=17"
$t.Cell(3,4).Range.Text = "42+44=86"
$t.Cell(3,5).Range.Text = "65+13=78"
$t.Cell(4,1).Range.Text = "99-57=42"
$t.Cell(4,2).Range.Text = "14-9=5"
$t.Cell(4,3).Range.Text = "28+6=34"
$t.Cell(4,4).Range.Text = "71-14=57"
$t.Cell(4,5).Range.Text = "97-79=18"
$t.Cell(5,1).Range.Text = "30+30=60"
$t.Cell(5,2).Range.Text = "55+41=96"
$t.Cell(5,3).Range.Text = "39+12=51"
$t.Cell(5,4).Range.Text = "65+29=94"
$t.Cell(5,5).Range.Text = "80-65=15"
$t.Cell(6,1).Range.Text = "55-36=19"
$t.Cell(6,2).Range.Text = "7+67=74"
$t.Cell(6,3).Range.Text = "0+38=38"
$t.Cell(6,4).Range.Text = "60-6=54"
$t.Cell(6,5).Range.Text = "46+39=85"
$t.Cell(7,1).Range.Text = "91+2=93"
$t.Cell(7,2).Range.Text = "61-12=49"
$t.Cell(7,3).Range.Text = "52+8=60"
$t.Cell(7,4).Range.Text = "18+48=66"
$t.Cell(7,5).Range.Text = "31+37=68"
$t.Cell(8,1).Range.Text = "64+25=89"
$t.Cell(8,2).Range.Text = "92-53=39"
$t.Cell(8,3).Range.Text = "64-9=55"
$t.Cell(8,4).Range.Text = "34+2=36"
$t.Cell(8,5).Range.Text = "98-61=37"
$t.Cell(9,1).Range.Text = "20+10=30"
$t.Cell(9,2).Range.Text = "63-26=37"
$t.Cell(9,3).Range.Text = "67-35=32"
$t.Cell(9,4).Range.Text = "11-4=7"
$t.Cell(9,5).Range.Text = "95-51=44"
$t.Cell(10,1).Range.Text = "36-12=24"
$t.Cell(10,2).Range.Text = "11+38=49"
$t.Cell(10,3).Range.Text = "84-16=68"
$t.Cell(10,4).Range.Text = "52-46=6"
$t.Cell(10,5).Range.Text = "73-19=54"
$t.Cell(11,1).Range.Text = "51-28=23"
$t.Cell(11,2).Range.Text = "75-11=64"
$t.Cell(11,3).Range.Text = "20+8=28"
$t.Cell(11,4).Range.Text = "75-2=73"
$t.Cell(11,5).Range.Text = "42-9=33"
$t.Cell(12,1).Range.Text = "19+34=53"
$t.Cell(12,2).Range.Text = "14+41=55"
$t.Cell(12,3).Range.Text = "50-26=24"
$t.Cell(12,4).Range.Text = "26+33=59"
$t.Cell(12,5).Range.Text = "65+34=99"
$t.Cell(13,1).Range.Text = "91-2=89"
$t.Cell(13,2).Range.Text = "28+68=96"
$t.Cell(13,3).Range.Text = "98-7=91"
$t.Cell(13,4).Range.Text = "21+48=69"
$t.Cell(13,5).Range.Text = "56+23=79"
$t.Cell(14,1).Range.Text = "32+42=74"
$t.Cell(14,2).Range.Text = "81-53=28"
$t.Cell(14,3).Range.Text = "23+22=45"
$t.Cell(14,4).Range.Text = "34+16=50"
$t.Cell(14,5).Range.Text = "16+56=72"
$t.Cell(15,1).Range.Text = "3+88=91"
$t.Cell(15,2).Range.Text = "39-5=34"
$t.Cell(15,3).Range.Text = "9+9=18"
$t.Cell(15,4).Range.Text = "37+2=39"
$t.Cell(15,5).Range.Text = "25+55=80"
$t.Cell(16,1).Range.Text = "25+30=55"
$t.Cell(16,2).Range.Text = "56+26=82"
$t.Cell(16,3).Range.Text = "78-17=61"
$t.Cell(16,4).Range.Text = "83-12=71"
$t.Cell(16,5).Range.Text = "22+32=54"
$t.Cell(17,1).Range.Text = "69+23=92"
$t.Cell(17,2).Range.Text = "23-7=16"
$t.Cell(17,3).Range.Text = "60-30=30"
$t.Cell(17,4).Range.Text = "18+74=92"
$t.Cell(17,5).Range.Text = "92-51=41"
$t.Cell(18,1).Range.Text = "72-48=24"
$t.Cell(18,2).Range.Text = "25+28=53"
$t.Cell(18,3).Range.Text = "91-10=81"
$t.Cell(18,4).Range.Text = "67+13=80"
$t.Cell(18,5).Range.Text = "86-33=53"
$t.Cell(19,1).Range.Text = "12+49=61"
$t.Cell(19,2).Range.Text = "2+10=12"
$t.Cell(19,3).Range.Text = "8+44=52"
$t.Cell(19,4).Range.Text = "44+24=68"
$t.Cell(19,5).Range.Text = "10-4=6"
$t.Cell(20,1).Range.Text = "46+53=99"
$t.Cell(20,2).Range.Text = "40+42=82"
$t.Cell(20,3).Range.Text = "33-13=20"
$t.Cell(20,4).Range.Text = "16+55=71"
$t.Cell(20,5).Range.Text = "37-3=34"
Write-Output "updated 100 cells"
